$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.174591510220848, 0.1258683508400295, 1.155256459293242, 0.1822841748126243),
    @(0.2954516743906898, 0.1361460328409026, 2.949074107730708, 0.29878681355207),
    @(0.2855239290104701, 0.1736533016683439, 3.666211514761394, 0.2230009963955268),
    @(0.1243648165791852, 0.114408573474128, 6.42049203935059, 0.1749400651923788),
    @(0.255523429346789, 0.1501314074596035, 8.770462758967463, 0.3770581559129562),
    @(0.2540495197352055, 0.09049827347454925, 9.727841100460605, 0.1928717669813325),
    @(0.1494001570811972, 0.1334653984824515, 12.46071583139859, 0.3121757356295127),
    @(0.2703745199553737, 0.1392095412808186, 14.69817297899096, 0.2863203349234984),
    @(0.1593764287736224, 0.09586511171837554, 15.28109735322164, 0.3658404188637184),
    @(0.2835086356192502, 0.09788018914535139, 17.54154182639366, 0.3052147677044286),
    @(0.1262559437430216, 0.1864472475251674, 19.37303394764616, 0.3943487833754562),
    @(0.1107120942924899, 0.0960975017415991, 20.73823033391738, 0.1751770002859284),
    @(0.2814512366711349, 0.1352677896355402, 21.69196349552623, 0.3385760466189176),
    @(0.2430667063229389, 0.1480316964732308, 25.69722107590501, 0.4142887397523345),
    @(0.2326655413894682, 0.1205392373751327, 27.0195648418906, 0.2819937502787892),
    @(0.1901363180535316, 0.09588772143902059, 29.58556960858154, 0.1834109107854891),
    @(0.1435525429874498, 0.1372984624797285, 30.19963124694031, 0.2613003229174928),
    @(0.1708694871536545, 0.1496943840011644, 32.85293938245817, 0.3886537178955867),
    @(0.2604398570439498, 0.2025188775800254, 35.18952747146144, 0.4313863417107523),
    @(0.1030433709033476, 0.1674002993442782, 35.62434345300911, 0.4426611090471849),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("B$r").Value = $data[$i][0]
    $ws.Range("D$r").Value = $data[$i][1]
    $ws.Range("E$r").Value = $data[$i][2]
    $ws.Range("F$r").Value = $data[$i][3]
}
